# remc_report_config.xlsx — "volt profile section data dump complete"
#
# 1. Rename the three existing section sheets to their new, more
#    descriptive names.
# 2. Append a new "volt_profile" worksheet after state_gen and populate it
#    with the 400kV / 220kV point data.
# 3. Re-point the active tab at ists_gen (the sheet that used to be
#    "section_2"), matching the saved workbook view.

$wb = $excel.ActiveWorkbook

$regional = $wb.Worksheets.Item(1)
$regional.Name = "regional_profile"

$ists = $wb.Worksheets.Item(2)
$ists.Name = "ists_gen"

$state = $wb.Worksheets.Item(3)
$state.Name = "state_gen"

# New sheet goes after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$volt = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$volt.Name = "volt_profile"

$volt.Range("A1").Value = "name"
$volt.Range("D1").Value = "type"
$volt.Range("B1").Value = "400_kv_pnt"
$volt.Range("C1").Value = "220_kv_pnt"

$volt.Range("D2").Value = "normal"
$volt.Range("B2").Value = "WREMCPRI.SCADA01.00016226"
$volt.Range("C2").Value = "WREMCPRI.SCADA01.00016228"
$volt.Range("A2").Value = "BACHHAU "

$volt.Range("D3").Value = "normal"
$volt.Range("B3").Value = "WREMCPRI.SCADA01.00016241"
$volt.Range("C3").Value = "WREMCPRI.SCADA01.00016243"
$volt.Range("A3").Value = "REWA"

$volt.Range("A4").Value = "Dummay Row"
$volt.Range("D4").Value = "dummy"

$volt.Columns.Item(1).AutoFit() | Out-Null
$volt.Columns.Item(2).AutoFit() | Out-Null
$volt.Columns.Item(3).AutoFit() | Out-Null

$volt.Range("D13").Select() | Out-Null

$state.Activate() | Out-Null
$state.Range("H26").Select() | Out-Null

# ists_gen ("section_2") is the tab that should end up active/selected.
$ists.Activate() | Out-Null
